# Remove the trailing footer block:
#   - the blank paragraph right after the last "Requisitos" line
#     ("LOQ4009: Instrumentação na Industria Química (Requisito fraco)")
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#      pages. Original theme under Creative Commons Attribution"
#
# After the edit, the "LOQ4009..." paragraph is immediately followed by
# the (still remaining) blank paragraph that precedes the page break.

$d = $word.ActiveDocument
$count = $d.Paragraphs.Count

# Locate the paragraph that starts the "LOQ4009..." requirement line.
$startIndex = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "LOQ4009*") {
        $startIndex = $i
        break
    }
}

# Locate the paragraph containing the copyright/footer text, which marks
# the end of the block to be removed.
$endIndex = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Creative Commons Attribution*") {
        $endIndex = $i
        break
    }
}

if ($startIndex -gt 0 -and $endIndex -ge $startIndex) {
    $firstToRemove = $d.Paragraphs.Item($startIndex + 1)
    $lastToRemove = $d.Paragraphs.Item($endIndex)
    $r = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
    $r.Delete()
}
